$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "_old" header suffixes (columns A-J) to "_FV2310"
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    $cell.Value = ($val -replace '_old$', '_FV2310')
}

# Rename the "_new" header suffixes (columns L-U) to "_FV2404"
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    $cell.Value = ($val -replace '_new$', '_FV2404')
}

# Turn the data range into a proper table with autofilter + banded rows
$rng = $ws.Range("A1:U85")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# Freeze the header row
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
